$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(9, 8).Value = 397.1905
$ws.Cells.Item(9, 9).Value = 221.72223
$ws.Cells.Item(9, 10).Value = 1450
$ws.Cells.Item(9, 11).Value = 221.72223
$ws.Cells.Item(9, 12).Value = 1450
$ws.Cells.Item(9, 13).Value = -52.72223
$ws.Cells.Item(9, 14).Value = -1788
$ws.Cells.Item(17, 8).Value = 1250.359
$ws.Cells.Item(17, 10).Value = 1366.6875
$ws.Cells.Item(17, 12).Value = 4100.0625
$ws.Cells.Item(17, 14).Value = -4436.0625
$ws.Cells.Item(19, 8).Value = 5360.609
$ws.Cells.Item(19, 9).Value = 5055.1113
$ws.Cells.Item(19, 10).Value = 5557
$ws.Cells.Item(19, 11).Value = 5055.1113
$ws.Cells.Item(19, 12).Value = 5557
$ws.Cells.Item(19, 13).Value = -4880.1113
$ws.Cells.Item(19, 14).Value = -5907
$ws.Cells.Item(33, 8).Value = 2412.0908
$ws.Cells.Item(33, 9).Value = 2862.7144
$ws.Cells.Item(33, 11).Value = 2862.7144
$ws.Cells.Item(33, 13).Value = -2633.7144
$ws.Cells.Item(112, 8).Value = 4961.5713
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 4961.5713
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 13).Value = 14884.7139
$ws.Cells.Item(112, 14).Value = -17100.7139
$ws.Cells.Item(135, 8).Value = 588981.25
$ws.Cells.Item(135, 9).Value = 714767.5600000001
$ws.Cells.Item(135, 11).Value = 6432908.040000001
$ws.Cells.Item(135, 13).Value = -6430373.040000001
$ws.Cells.Item(137, 8).Value = 3114.8774
$ws.Cells.Item(137, 9).Value = 2686.3462
$ws.Cells.Item(137, 11).Value = 8059.0386
$ws.Cells.Item(137, 13).Value = -5509.0386
$ws.Cells.Item(138, 8).Value = 1566603.4
$ws.Cells.Item(138, 10).Value = 2784080
$ws.Cells.Item(138, 12).Value = 8352240
$ws.Cells.Item(138, 14).Value = -8362520

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(16, 8).Value = 4423.75
$ws.Cells.Item(16, 9).Value = 1196.6666
$ws.Cells.Item(16, 10).Value = 6360
$ws.Cells.Item(16, 11).Value = 1196.6666
$ws.Cells.Item(16, 12).Value = 6360
$ws.Cells.Item(16, 13).Value = -909.6666
$ws.Cells.Item(16, 14).Value = -6934
$ws.Cells.Item(32, 8).Value = 2237863
$ws.Cells.Item(32, 9).Value = 2319357.8
$ws.Cells.Item(32, 11).Value = 2319357.8
$ws.Cells.Item(32, 13).Value = -2319070.8
$ws.Cells.Item(61, 8).Value = 5900.1665
$ws.Cells.Item(61, 9).Value = 3738
$ws.Cells.Item(61, 11).Value = 3738
$ws.Cells.Item(61, 13).Value = -3526
$ws.Cells.Item(132, 8).Value = 4096.028
$ws.Cells.Item(132, 9).Value = 1733.1489
$ws.Cells.Item(132, 11).Value = 5199.4467
$ws.Cells.Item(132, 13).Value = -2669.4467
$ws.Cells.Item(136, 8).Value = 5900.1665
$ws.Cells.Item(136, 9).Value = 3738
$ws.Cells.Item(136, 11).Value = 11214
$ws.Cells.Item(136, 13).Value = -8664

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 4775.1333
$ws.Cells.Item(134, 9).Value = 2130.6287
$ws.Cells.Item(134, 10).Value = 8477.440000000001
$ws.Cells.Item(134, 11).Value = 6391.886100000001
$ws.Cells.Item(134, 12).Value = 25432.32
$ws.Cells.Item(134, 13).Value = -3856.886100000001
$ws.Cells.Item(134, 14).Value = -30502.32

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(14, 8).Value = 565
$ws.Cells.Item(22, 8).Value = 586
$ws.Cells.Item(22, 9).Value = 593.75
$ws.Cells.Item(22, 10).Value = 579.8
$ws.Cells.Item(22, 11).Value = 593.75
$ws.Cells.Item(22, 12).Value = 579.8
$ws.Cells.Item(22, 13).Value = -243.75
$ws.Cells.Item(22, 14).Value = -1279.8
$ws.Cells.Item(31, 8).Value = 5959380
$ws.Cells.Item(31, 9).Value = 2963.2903
$ws.Cells.Item(31, 10).Value = 13345337
$ws.Cells.Item(31, 11).Value = 2963.2903
$ws.Cells.Item(31, 12).Value = 13345337
$ws.Cells.Item(31, 13).Value = -2668.2903
$ws.Cells.Item(31, 14).Value = -13345927
$ws.Cells.Item(34, 8).Value = 5959380
$ws.Cells.Item(34, 9).Value = 2963.2903
$ws.Cells.Item(34, 10).Value = 13345337
$ws.Cells.Item(34, 11).Value = 2963.2903
$ws.Cells.Item(34, 12).Value = 13345337
$ws.Cells.Item(34, 13).Value = -2761.2903
$ws.Cells.Item(34, 14).Value = -13345741
$ws.Cells.Item(134, 8).Value = 3692.5764
$ws.Cells.Item(134, 9).Value = 1403.138
$ws.Cells.Item(134, 10).Value = 8610.629999999999
$ws.Cells.Item(134, 11).Value = 4209.414
$ws.Cells.Item(134, 12).Value = 25831.89
$ws.Cells.Item(134, 13).Value = -1674.414
$ws.Cells.Item(134, 14).Value = -30901.89

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(87, 13).ClearContents()
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(5, 8).Value = 5028.7
$ws.Cells.Item(5, 9).Value = 2497
$ws.Cells.Item(5, 10).Value = 6716.5
$ws.Cells.Item(5, 11).Value = 7491
$ws.Cells.Item(5, 12).Value = 20149.5
$ws.Cells.Item(5, 13).Value = -7379
$ws.Cells.Item(5, 14).Value = -20373.5
$ws.Cells.Item(12, 8).Value = 2381620.8
$ws.Cells.Item(12, 9).Value = 1903
$ws.Cells.Item(12, 10).Value = 3333508
$ws.Cells.Item(12, 11).Value = 5709
$ws.Cells.Item(12, 12).Value = 10000524
$ws.Cells.Item(12, 13).Value = -5536
$ws.Cells.Item(12, 14).Value = -10000870
$ws.Cells.Item(14, 8).Value = 18519926
$ws.Cells.Item(14, 9).Value = 18519926
$ws.Cells.Item(14, 11).Value = 55559778
$ws.Cells.Item(14, 13).Value = -55559605
$ws.Cells.Item(38, 8).Value = 50000036
$ws.Cells.Item(38, 10).Value = 125000024
$ws.Cells.Item(38, 12).Value = 375000072
$ws.Cells.Item(38, 14).Value = -375000766
$ws.Cells.Item(61, 8).Value = 357.5
$ws.Cells.Item(61, 9).Value = 92.5
$ws.Cells.Item(61, 10).Value = 887.5
$ws.Cells.Item(61, 11).Value = 277.5
$ws.Cells.Item(61, 12).Value = 2662.5
$ws.Cells.Item(61, 13).Value = -62.5
$ws.Cells.Item(61, 14).Value = -3092.5
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(92, 8).Value = 6411881
$ws.Cells.Item(92, 9).Value = 1000
$ws.Cells.Item(92, 10).Value = 6994688
$ws.Cells.Item(92, 11).Value = 3000
$ws.Cells.Item(92, 12).Value = 20984064
$ws.Cells.Item(92, 13).Value = -1752
$ws.Cells.Item(92, 14).Value = -20986560
$ws.Cells.Item(97, 8).Value = 362.6
$ws.Cells.Item(97, 9).Value = 362.6
$ws.Cells.Item(97, 11).Value = 1087.8
$ws.Cells.Item(97, 13).Value = -591.8000000000002
$ws.Cells.Item(122, 8).Value = 1230864.2
$ws.Cells.Item(122, 9).Value = 2176842.8
$ws.Cells.Item(122, 10).Value = 1092
$ws.Cells.Item(122, 11).Value = 19591585.2
$ws.Cells.Item(122, 12).Value = 9828
$ws.Cells.Item(122, 13).Value = -19589135.2
$ws.Cells.Item(122, 14).Value = -14728
$ws.Cells.Item(131, 8).Value = 2209.9524
$ws.Cells.Item(131, 10).Value = 2291.3242
$ws.Cells.Item(131, 12).Value = 6873.9726
$ws.Cells.Item(131, 14).Value = -16953.9726
$ws.Cells.Item(135, 8).Value = 5028.7
$ws.Cells.Item(135, 9).Value = 2497
$ws.Cells.Item(135, 10).Value = 6716.5
$ws.Cells.Item(135, 11).Value = 22473
$ws.Cells.Item(135, 12).Value = 60448.5
$ws.Cells.Item(135, 13).Value = -19938
$ws.Cells.Item(135, 14).Value = -65518.5

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 1176698.6
$ws.Cells.Item(2, 9).Value = 69.85714
$ws.Cells.Item(2, 10).Value = 2000338.8
$ws.Cells.Item(2, 11).Value = 69.85714
$ws.Cells.Item(2, 12).Value = 2000338.8
$ws.Cells.Item(2, 13).Value = 43.14286
$ws.Cells.Item(2, 14).Value = -2000564.8
$ws.Cells.Item(33, 8).Value = 30666.666
$ws.Cells.Item(33, 10).Value = 30666.666
$ws.Cells.Item(33, 12).Value = 30666.666
$ws.Cells.Item(33, 14).Value = -31170.666
$ws.Cells.Item(35, 8).Value = 29300
$ws.Cells.Item(35, 10).Value = 29300
$ws.Cells.Item(35, 12).Value = 29300
$ws.Cells.Item(35, 14).Value = -29896
$ws.Cells.Item(58, 8).Value = 58489.668
$ws.Cells.Item(58, 10).Value = 66379.60000000001
$ws.Cells.Item(58, 12).Value = 66379.60000000001
$ws.Cells.Item(58, 14).Value = -66933.60000000001
$ws.Cells.Item(102, 8).Value = 1841.5319
$ws.Cells.Item(102, 9).Value = 1487.85
$ws.Cells.Item(102, 10).Value = 3862.5715
$ws.Cells.Item(102, 11).Value = 1487.85
$ws.Cells.Item(102, 12).Value = 3862.5715
$ws.Cells.Item(102, 13).Value = 134.1500000000001
$ws.Cells.Item(102, 14).Value = -7106.5715

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 5013.189
$ws.Cells.Item(40, 9).Value = 3473.2173
$ws.Cells.Item(40, 11).Value = 3473.2173
$ws.Cells.Item(40, 13).Value = -3337.2173
$ws.Cells.Item(68, 8).Value = 4530.5386
$ws.Cells.Item(68, 10).Value = 4670.9414
$ws.Cells.Item(68, 12).Value = 4670.9414
$ws.Cells.Item(68, 14).Value = -6168.9414
$ws.Cells.Item(71, 8).Value = 4530.5386
$ws.Cells.Item(71, 10).Value = 4670.9414
$ws.Cells.Item(71, 12).Value = 23354.707
$ws.Cells.Item(71, 14).Value = -30842.707
$ws.Cells.Item(100, 8).Value = 3516
$ws.Cells.Item(100, 9).Value = 2650.4285
$ws.Cells.Item(100, 11).Value = 2650.4285
$ws.Cells.Item(100, 13).Value = -2109.4285
$ws.Cells.Item(132, 8).Value = 9440186
$ws.Cells.Item(132, 9).Value = 21741714
$ws.Cells.Item(132, 11).Value = 65225142
$ws.Cells.Item(132, 13).Value = -65222612
$ws.Cells.Item(136, 8).Value = 5771.0493
$ws.Cells.Item(136, 9).Value = 1857.1945
$ws.Cells.Item(136, 11).Value = 5571.583500000001
$ws.Cells.Item(136, 13).Value = -3021.583500000001

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 30772244
$ws.Cells.Item(132, 9).Value = 35090496
$ws.Cells.Item(132, 11).Value = 105271488
$ws.Cells.Item(132, 13).Value = -105268958
$ws.Cells.Item(136, 8).Value = 16148212
$ws.Cells.Item(136, 9).Value = 24390956
$ws.Cells.Item(136, 11).Value = 73172868
$ws.Cells.Item(136, 13).Value = -73170318
